# Updates cryptocurrency price/volume figures (and, for rows 43-45, the coin
# name/link) on the single worksheet of the active workbook, matching the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell -> new text value. Column D (Price) entries that look like a
# plain number (e.g. "1.00", "0.880") would otherwise be auto-converted by
# Excel into a numeric value and lose formatting such as trailing zeros, so
# those cells are temporarily switched to Text format, written, and then
# reset back to the default "Normal" style so no visible formatting changes.
$cellUpdates = [ordered]@{
    "D2" = "64.523.57"
    "E2" = "  +0.99%  "
    "D3" = "3.169.52"
    "E4" = "  -0.03%  "
    "D5" = "572.53"
    "E5" = "  +1.12%  "
    "D6" = "164.42"
    "E6" = "  -1.88%  "
    "E7" = "  +0.02%  "
    "E8" = "  -3.62%  "
    "E9" = "  -2.05%  "
    "D10" = "6.62"
    "E10" = "  -0.89%  "
    "E11" = "  +0.07%  "
    "D12" = "3.724.14"
    "E12" = "  +0.13%  "
    "E13" = "  -1.02%  "
    "D14" = "64.545.98"
    "E14" = "  +0.87%  "
    "D15" = "25.33"
    "E15" = "  +0.23%  "
    "D16" = "3.171.44"
    "E16" = "  +0.23%  "
    "E17" = "  -1.71%  "
    "D18" = "407.34"
    "E19" = "  +0.01%  "
    "D20" = "5.25"
    "E20" = "  -1.24%  "
    "D21" = "7.15"
    "E21" = "  +0.80%  "
    "E22" = "  +0.07%  "
    "D23" = "68.77"
    "E23" = "  -2.83%  "
    "D24" = "0.487"
    "E24" = "  -0.51%  "
    "D25" = "0.194"
    "E25" = "  -3.68%  "
    "E26" = "  -3.71%  "
    "D27" = "8.85"
    "D28" = "1.00"
    "E28" = "  +0.10%  "
    "E29" = "  -0.25%  "
    "D30" = "21.28"
    "E30" = "  -2.09%  "
    "D31" = "6.34"
    "E31" = "  +0.23%  "
    "E32" = "  -1.89%  "
    "E33" = "  +0.24%  "
    "D34" = "156.55"
    "E34" = "  +0.56%  "
    "E35" = "  -1.12%  "
    "E36" = "  +0.39%  "
    "D37" = "2.691.88"
    "E37" = "  -1.26%  "
    "D38" = "24.07"
    "E38" = "  -2.19%  "
    "E39" = "  -1.01%  "
    "E40" = "  -2.07%  "
    "D41" = "0.0619"
    "E41" = "  -0.19%  "
    "D42" = "5.52"
    "E42" = "  -1.69%  "
    "B43" = "Bittensor"
    "C43" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D43" = "291.80"
    "E43" = "  -0.68%  "
    "B44" = "InjectiveProtocol"
    "C44" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D44" = "21.49"
    "E44" = "  -1.29%  "
    "B45" = "VeChain"
    "C45" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D45" = "0.0258"
    "E45" = "  -1.38%  "
    "E46" = "  +0.00%  "
    "D47" = "0.0984"
    "E47" = "  -0.26%  "
    "E48" = "  -5.49%  "
    "D49" = "10.45"
    "E49" = "  +0.34%  "
    "E50" = "  -0.66%  "
    "D51" = "0.880"
    "E51" = "  -6.16%  "
}

foreach ($cellRef in $cellUpdates.Keys) {
    $newValue = $cellUpdates[$cellRef]
    $range = $ws.Range($cellRef)

    $needsTextFormat = ($cellRef.StartsWith("D")) -and ($newValue -match "^\s*[-+]?[0-9]*\.?[0-9]+\s*$")
    if ($needsTextFormat) {
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
